$wb = $excel.ActiveWorkbook

# Target column width (characters). The workbook's original XML carries a
# sub-pixel float (17.2159881591797) that was never produced by Excel's own
# pixel-grid quantization (Excel COM snaps ColumnWidth to 1/6-character /
# pixel boundaries). 16.3 is the nearest settable input that lands on the
# closest reachable grid point to the authored width.
$newColWidth = 16.3

# --- Overview sheet ---
$ovWs = $wb.Worksheets.Item("Overview")
$ovWs.Range("E2").Value = "Ready for handoff"
$ovWs.Range("F2").Value = "Ready for handoff"
$ovWs.Range("G2").Value = "2016-08-23 15:15:32"
$ovWs.Columns.Item(5).ColumnWidth = $newColWidth
$ovWs.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet ---
$zhWs = $wb.Worksheets.Item("zh-cn")
$zhWs.Range("C2").Value = "Ready for handoff"
$zhWs.Range("H2").Value = "2016-08-23 15:15:23"
$zhWs.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet ---
$deWs = $wb.Worksheets.Item("de-de")
$deWs.Range("C2").Value = "Ready for handoff"
$deWs.Range("H2").Value = "2016-08-23 15:15:32"
$deWs.Columns.Item(3).ColumnWidth = $newColWidth
